$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.263122
$ws.Range("H2").Value = 9.789365999999999
$ws.Range("I2").Value = 0.3531375780718168
$ws.Range("J2").Value = 0.3531375780718168
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 35.31114333333333
$ws.Range("N2").Value = 105.93343
$ws.Range("O2").Value = 0.6187867769880316
$ws.Range("P2").Value = 0.6187867769880316
$ws.Range("Q2").Value = 115.2245686561533
$ws.Range("R2").Value = 1037.02111790538
$ws.Range("S2").Value = 0.2185168637684189
$ws.Range("T2").Value = 0.2185168637684189

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.263122
$ws.Range("H3").Value = 9.789365999999999
$ws.Range("I3").Value = 0.3531375780718168
$ws.Range("J3").Value = 0.3531375780718168
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.551362
$ws.Range("N3").Value = 19.654086
$ws.Range("O3").Value = 0.1148050103785518
$ws.Range("P3").Value = 0.1148050103785518
$ws.Range("Q3").Value = 21.377893472164
$ws.Range("R3").Value = 192.401041249476
$ws.Range("S3").Value = 0.04054196331559158
$ws.Range("T3").Value = 0.04054196331559158

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.263122
$ws.Range("H4").Value = 9.789365999999999
$ws.Range("I4").Value = 0.3531375780718168
$ws.Range("J4").Value = 0.3531375780718168
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.951915
$ws.Range("N4").Value = 32.855745
$ws.Range("O4").Value = 0.191919590955288
$ws.Range("P4").Value = 0.191919590955288
$ws.Range("Q4").Value = 35.73743477863
$ws.Range("R4").Value = 321.63691300767
$ws.Range("S4").Value = 0.06777401953448416
$ws.Range("T4").Value = 0.06777401953448416

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.263122
$ws.Range("H5").Value = 9.789365999999999
$ws.Range("I5").Value = 0.3531375780718168
$ws.Range("J5").Value = 0.3531375780718168
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.250702333333333
$ws.Range("N5").Value = 12.752107
$ws.Range("O5").Value = 0.07448862167812857
$ws.Range("P5").Value = 0.07448862167812857
$ws.Range("Q5").Value = 13.87056029935133
$ws.Range("R5").Value = 124.835042694162
$ws.Range("S5").Value = 0.02630473145332215
$ws.Range("T5").Value = 0.02630473145332215

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.367310666666667
$ws.Range("H6").Value = 13.101932
$ws.Range("I6").Value = 0.4726337266929886
$ws.Range("J6").Value = 0.4726337266929886
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.31114333333333
$ws.Range("N6").Value = 105.93343
$ws.Range("O6").Value = 0.6187867769880316
$ws.Range("P6").Value = 0.6187867769880316
$ws.Range("Q6").Value = 154.2147329318622
$ws.Range("R6").Value = 1387.93259638676
$ws.Range("S6").Value = 0.2924595004361966
$ws.Range("T6").Value = 0.2924595004361966

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.367310666666667
$ws.Range("H7").Value = 13.101932
$ws.Range("I7").Value = 0.4726337266929886
$ws.Range("J7").Value = 0.4726337266929886
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.551362
$ws.Range("N7").Value = 19.654086
$ws.Range("O7").Value = 0.1148050103785518
$ws.Range("P7").Value = 0.1148050103785518
$ws.Range("Q7").Value = 28.61183314379467
$ws.Range("R7").Value = 257.506498294152
$ws.Range("S7").Value = 0.05426071989824218
$ws.Range("T7").Value = 0.05426071989824218

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.367310666666667
$ws.Range("H8").Value = 13.101932
$ws.Range("I8").Value = 0.4726337266929886
$ws.Range("J8").Value = 0.4726337266929886
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.951915
$ws.Range("N8").Value = 32.855745
$ws.Range("O8").Value = 0.191919590955288
$ws.Range("P8").Value = 0.191919590955288
$ws.Range("Q8").Value = 47.83041519992666
$ws.Range("R8").Value = 430.47373679934
$ws.Range("S8").Value = 0.09070767149859174
$ws.Range("T8").Value = 0.09070767149859174

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.367310666666667
$ws.Range("H9").Value = 13.101932
$ws.Range("I9").Value = 0.4726337266929886
$ws.Range("J9").Value = 0.4726337266929886
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.250702333333333
$ws.Range("N9").Value = 12.752107
$ws.Range("O9").Value = 0.07448862167812857
$ws.Range("P9").Value = 0.07448862167812857
$ws.Range("Q9").Value = 18.56413764119155
$ws.Range("R9").Value = 167.077238770724
$ws.Range("S9").Value = 0.03520583485995805
$ws.Range("T9").Value = 0.03520583485995805

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.609937666666666
$ws.Range("H10").Value = 4.829813
$ws.Range("I10").Value = 0.1742286952351946
$ws.Range("J10").Value = 0.1742286952351946
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 35.31114333333333
$ws.Range("N10").Value = 105.93343
$ws.Range("O10").Value = 0.6187867769880316
$ws.Range("P10").Value = 0.6187867769880316
$ws.Range("Q10").Value = 56.84873970539888
$ws.Range("R10").Value = 511.63865734859
$ws.Range("S10").Value = 0.1078104127834161
$ws.Range("T10").Value = 0.1078104127834161

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.609937666666666
$ws.Range("H11").Value = 4.829813
$ws.Range("I11").Value = 0.1742286952351946
$ws.Range("J11").Value = 0.1742286952351946
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.551362
$ws.Range("N11").Value = 19.654086
$ws.Range("O11").Value = 0.1148050103785518
$ws.Range("P11").Value = 0.1148050103785518
$ws.Range("Q11").Value = 10.54728445176867
$ws.Range("R11").Value = 94.92556006591799
$ws.Range("S11").Value = 0.02000232716471806
$ws.Range("T11").Value = 0.02000232716471806

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.609937666666666
$ws.Range("H12").Value = 4.829813
$ws.Range("I12").Value = 0.1742286952351946
$ws.Range("J12").Value = 0.1742286952351946
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.951915
$ws.Range("N12").Value = 32.855745
$ws.Range("O12").Value = 0.191919590955288
$ws.Range("P12").Value = 0.191919590955288
$ws.Range("Q12").Value = 17.63190048063166
$ws.Range("R12").Value = 158.687104325685
$ws.Range("S12").Value = 0.03343789992221208
$ws.Range("T12").Value = 0.03343789992221208

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.609937666666666
$ws.Range("H13").Value = 4.829813
$ws.Range("I13").Value = 0.1742286952351946
$ws.Range("J13").Value = 0.1742286952351946
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.250702333333333
$ws.Range("N13").Value = 12.752107
$ws.Range("O13").Value = 0.07448862167812857
$ws.Range("P13").Value = 0.07448862167812857
$ws.Range("Q13").Value = 6.843365796221221
$ws.Range("R13").Value = 61.59029216599099
$ws.Range("S13").Value = 0.01297805536484837
$ws.Range("T13").Value = 0.01297805536484837
